# Atualização automática: 2025-09-02 21:00:26
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 7-11 rotate one step: row 7's original record moves down to row 11,
# and the records that were in rows 8-11 each shift up one row (8->7, 9->8,
# 10->9, 11->10). Columns B (Class) and C (First_Detection_Date) are
# unaffected by the rotation.

$cols = @("A", "D", "E", "F", "G", "H", "I", "J")

# Make sure the coords (I) and confidence (J) columns keep their original
# text type instead of being auto-coerced to numbers by Excel.
$ws.Range("I7:J11").NumberFormat = "@"
$ws.Range("I18:J18").NumberFormat = "@"

# Capture the original row 7 values before overwriting anything.
$orig7 = @{}
foreach ($col in $cols) {
    $addr = "${col}7"
    $orig7[$col] = $ws.Range($addr).Value()
}

# Shift rows 8..11 up into rows 7..10.
for ($r = 7; $r -le 10; $r++) {
    $srcRow = $r + 1
    foreach ($col in $cols) {
        $dstAddr = "${col}${r}"
        $srcAddr = "${col}${srcRow}"
        $ws.Range($dstAddr).Value = $ws.Range($srcAddr).Value()
    }
}

# Place the original row 7 values into row 11.
foreach ($col in $cols) {
    $addr = "${col}11"
    $ws.Range($addr).Value = $orig7[$col]
}

# Row 18: update the detection image, bounding-box coords, and confidence.
$ws.Range("D18").Value = "image_20250808100711_ppp0.jpg"
$ws.Range("I18").Value = "1182,409,1232,451"
$ws.Range("J18").Value = "0.75"
